$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E3").Formula = "=(H11+J11)*E11/(PI())"
$ws.Range("F3").Formula = "=(I11+K11)*E11/PI()"

$ws.Range("E4").Formula = "=(H12+J12)*E12/(PI())"
$ws.Range("F4").Formula = "=(I12+K12)*E12/PI()"

$ws.Range("E5").Formula = "=(H13+J13)*E13/(PI())"
$ws.Range("F5").Formula = "=(I13+K13)*E13/PI()"

$ws.Range("E6").Formula = "=(H14+J14)*E14/(PI())"
$ws.Range("F6").Formula = "=(I14+K14)*E14/PI()"

$ws.Range("E7").Formula = "=(H15+J15)*E15/(PI())"
$ws.Range("F7").Formula = "=(I15+K15)*E15/PI()"
